$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 15-16 (existing rows 15.. shift down to 17..)
$ws.Rows("15:16").Insert()

# Seed the new rows with the same constant columns/formatting as the (now) following
# data rows 17:18, which used to be rows 15:16 before the insert.
$ws.Range("A17:T18").Copy()
$ws.Range("A15").PasteSpecial()

# Now overwrite the week-specific figures for the two new rows (new weekly price
# observations for Kiwi, Primera/Segunda).
# Row 15 - Primera
$ws.Range("D15").Value = 44414
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 12500
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 12750
$ws.Range("S15").Value = 708

# Row 16 - Segunda
$ws.Range("D16").Value = 44414
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11500
$ws.Range("P16").Value = 11250
$ws.Range("S16").Value = 625
